$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All updated cells in this sheet are stored as text (inlineStr) in the
# original workbook. Column D in particular holds values that LOOK numeric
# (e.g. "1.005", "0.00001009") but must stay as literal text. Briefly force
# the cell to a text number format while assigning, then restore the
# original "General" format so no stray formatting diff is introduced.

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '20.413.06'
$ws.Range('D2').NumberFormat = 'General'
$ws.Range('E2').Value = '  -7.30%  '
# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.443.30'
$ws.Range('D3').NumberFormat = 'General'
$ws.Range('E3').Value = '  -7.02%  '
# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.005'
$ws.Range('D4').NumberFormat = 'General'
$ws.Range('E4').Value = '  +0.29%  '
# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '1.004'
$ws.Range('D5').NumberFormat = 'General'
$ws.Range('E5').Value = '  +0.24%  '
# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '277.94'
$ws.Range('D6').NumberFormat = 'General'
$ws.Range('E6').Value = '  -4.03%  '
# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3746'
$ws.Range('D7').NumberFormat = 'General'
$ws.Range('E7').Value = '  -4.64%  '
# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3072'
$ws.Range('D8').NumberFormat = 'General'
$ws.Range('E8').Value = '  -4.07%  '
# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '40.61'
$ws.Range('D9').NumberFormat = 'General'
$ws.Range('E9').Value = '  -8.34%  '
# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.013'
$ws.Range('D10').NumberFormat = 'General'
$ws.Range('E10').Value = '  -5.22%  '
# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.06551'
$ws.Range('D11').NumberFormat = 'General'
$ws.Range('E11').Value = '  -8.75%  '
# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.005'
$ws.Range('D12').NumberFormat = 'General'
$ws.Range('E12').Value = '  +0.29%  '
# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.380'
$ws.Range('D13').NumberFormat = 'General'
$ws.Range('E13').Value = '  -4.53%  '
# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '17.26'
$ws.Range('D14').NumberFormat = 'General'
$ws.Range('E14').Value = '  -6.98%  '
# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.137'
$ws.Range('D15').NumberFormat = 'General'
$ws.Range('E15').Value = '  -7.51%  '
# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.444.44'
$ws.Range('D16').NumberFormat = 'General'
$ws.Range('E16').Value = '  -7.10%  '
# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001009'
$ws.Range('D17').NumberFormat = 'General'
$ws.Range('E17').Value = '  -7.76%  '
# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.05870'
$ws.Range('D18').NumberFormat = 'General'
$ws.Range('E18').Value = '  -10.49%  '
# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '76.15'
$ws.Range('D19').NumberFormat = 'General'
$ws.Range('E19').Value = '  -8.45%  '
# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.004'
$ws.Range('D20').NumberFormat = 'General'
$ws.Range('E20').Value = '  +0.36%  '
# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.734'
$ws.Range('D21').NumberFormat = 'General'
$ws.Range('E21').Value = '  -7.10%  '
# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '14.39'
$ws.Range('D22').NumberFormat = 'General'
$ws.Range('E22').Value = '  -6.44%  '
# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.89'
$ws.Range('D23').NumberFormat = 'General'
$ws.Range('E23').Value = '  -2.13%  '
# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.310'
$ws.Range('D24').NumberFormat = 'General'
$ws.Range('E24').Value = '  -2.42%  '
# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '20.415.12'
$ws.Range('D25').NumberFormat = 'General'
$ws.Range('E25').Value = '  -7.33%  '
# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '143.61'
$ws.Range('D26').NumberFormat = 'General'
$ws.Range('E26').Value = '  -2.64%  '
# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.217'
$ws.Range('D27').NumberFormat = 'General'
$ws.Range('E27').Value = '  -6.68%  '
# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '17.02'
$ws.Range('D28').NumberFormat = 'General'
$ws.Range('E28').Value = '  -7.61%  '
# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.607.37'
$ws.Range('D29').NumberFormat = 'General'
$ws.Range('E29').Value = '  -7.04%  '
# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '109.46'
$ws.Range('D30').NumberFormat = 'General'
$ws.Range('E30').Value = '  -7.13%  '
# Row 31
$ws.Range('B31').Value = 'HuobiToken'
$ws.Range('C31').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.756'
$ws.Range('D31').NumberFormat = 'General'
$ws.Range('E31').Value = '  -22.85%  '
# Row 32
$ws.Range('B32').Value = 'ImmutableX'
$ws.Range('C32').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.9028'
$ws.Range('D32').NumberFormat = 'General'
$ws.Range('E32').Value = '  -7.11%  '
# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.421'
$ws.Range('D33').NumberFormat = 'General'
$ws.Range('E33').Value = '  -6.34%  '
# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.07741'
$ws.Range('D34').NumberFormat = 'General'
$ws.Range('E34').Value = '  -6.61%  '
# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '8.279'
$ws.Range('D35').NumberFormat = 'General'
$ws.Range('E35').Value = '  -9.00%  '
# Row 36
$ws.Range('B36').Value = 'Frax'
$ws.Range('C36').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.004'
$ws.Range('D36').NumberFormat = 'General'
$ws.Range('E36').Value = '  +0.32%  '
# Row 37
$ws.Range('B37').Value = 'Aptos'
$ws.Range('C37').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '10.87'
$ws.Range('D37').NumberFormat = 'General'
$ws.Range('E37').Value = '  +2.06%  '
# Row 38
$ws.Range('B38').Value = 'Hedera'
$ws.Range('C38').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.05626'
$ws.Range('D38').NumberFormat = 'General'
$ws.Range('E38').Value = '  -5.90%  '
# Row 39
$ws.Range('B39').Value = 'TrustWalletToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.144'
$ws.Range('D39').NumberFormat = 'General'
$ws.Range('E39').Value = '  -5.44%  '
# Row 40
$ws.Range('B40').Value = 'InternetComputer(DFINITY)'
$ws.Range('C40').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '4.733'
$ws.Range('D40').NumberFormat = 'General'
$ws.Range('E40').Value = '  -6.87%  '
# Row 41
$ws.Range('B41').Value = 'WEMIXTOKEN'
$ws.Range('C41').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.375'
$ws.Range('D41').NumberFormat = 'General'
$ws.Range('E41').Value = '  -14.30%  '
# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1915'
$ws.Range('D42').NumberFormat = 'General'
$ws.Range('E42').Value = '  -6.34%  '
# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.02038'
$ws.Range('D43').NumberFormat = 'General'
$ws.Range('E43').Value = '  -9.07%  '
# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.590'
$ws.Range('D44').NumberFormat = 'General'
$ws.Range('E44').Value = '  -4.21%  '
# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.5318'
$ws.Range('D45').NumberFormat = 'General'
$ws.Range('E45').Value = '  -7.95%  '
# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '12.08'
$ws.Range('D46').NumberFormat = 'General'
$ws.Range('E46').Value = '  -6.98%  '
# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5155'
$ws.Range('D47').NumberFormat = 'General'
$ws.Range('E47').Value = '  -6.80%  '
# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '111.55'
$ws.Range('D48').NumberFormat = 'General'
$ws.Range('E48').Value = '  -4.52%  '
# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.789'
$ws.Range('D49').NumberFormat = 'General'
$ws.Range('E49').Value = '  -4.11%  '
# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.053'
$ws.Range('D50').NumberFormat = 'General'
$ws.Range('E50').Value = '  -6.75%  '
# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.004'
$ws.Range('D51').NumberFormat = 'General'
$ws.Range('E51').Value = '  +0.31%  '
